$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-18"

# Update the header label for the current-year column (I1)
$ws.Range("I1").Value = "2022 (through 06-18)"

# Update the June (row 7) value for the current-year column
$ws.Range("I7").Value = 92

# Update the Total (row 14) value for the current-year column
$ws.Range("I14").Value = 755
